# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts and one "最低票价" (G column)
# price correction across the "展览", "本地生活" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 506
$ws.Range("F7").Value = 9621
$ws.Range("F9").Value = 721
$ws.Range("F10").Value = 2188
$ws.Range("F12").Value = 1656
$ws.Range("F13").Value = 2781
$ws.Range("F15").Value = 4155
$ws.Range("F19").Value = 527
$ws.Range("F20").Value = 246
$ws.Range("F21").Value = 39
$ws.Range("F23").Value = 88
$ws.Range("F25").Value = 4026
$ws.Range("F27").Value = 3487
$ws.Range("G27").Value = 70
$ws.Range("F28").Value = 1110
$ws.Range("F29").Value = 203
$ws.Range("F30").Value = 509
$ws.Range("F31").Value = 4351
$ws.Range("F33").Value = 349
$ws.Range("F34").Value = 437
$ws.Range("F35").Value = 333

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 203
$ws.Range("F3").Value = 1012

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 203
$ws.Range("F4").Value = 1012
$ws.Range("F7").Value = 506
$ws.Range("F9").Value = 9621
$ws.Range("F11").Value = 721
$ws.Range("F12").Value = 2188
$ws.Range("F14").Value = 1656
$ws.Range("F16").Value = 2781
$ws.Range("F18").Value = 4155
$ws.Range("F22").Value = 527
$ws.Range("F23").Value = 246
$ws.Range("F24").Value = 39
$ws.Range("F27").Value = 88
$ws.Range("F29").Value = 4026
$ws.Range("F31").Value = 3487
$ws.Range("G31").Value = 70
$ws.Range("F32").Value = 1110
$ws.Range("F33").Value = 203
$ws.Range("F34").Value = 509
$ws.Range("F35").Value = 4351
$ws.Range("F37").Value = 349
$ws.Range("F38").Value = 437
$ws.Range("F39").Value = 333
